$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as Text so that numeric-looking
# values (e.g. "8.42") are stored as strings, matching the source data
# which is exported as literal text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.988.09'
$ws.Range("E2").Value = '  +2.79%  '

$ws.Range("D3").Value = '3.293.07'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '582.74'
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("D6").Value = '181.92'
$ws.Range("E6").Value = '  -2.13%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +2.13%  '

$ws.Range("D9").Value = '3.289.78'
$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("D11").Value = '0.579'
$ws.Range("E11").Value = '  +0.77%  '

$ws.Range("D12").Value = '46.03'
$ws.Range("E12").Value = '  -1.02%  '

$ws.Range("E13").Value = '  +2.95%  '

$ws.Range("D14").Value = '680.11'
$ws.Range("E14").Value = '  +10.64%  '

$ws.Range("D15").Value = '3.826.20'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = '8.42'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").Value = '68.087.98'
$ws.Range("E17").Value = '  +2.90%  '

$ws.Range("E18").Value = '  +1.36%  '

$ws.Range("D19").Value = '3.304.13'
$ws.Range("E19").Value = '  -0.53%  '

$ws.Range("D20").Value = '17.52'
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("D21").Value = '10.91'
$ws.Range("E21").Value = '  -0.61%  '

$ws.Range("D22").Value = '0.897'
$ws.Range("E22").Value = '  +0.23%  '

$ws.Range("D23").Value = '17.58'
$ws.Range("E23").Value = '  -3.43%  '

$ws.Range("D24").Value = '5.16'
$ws.Range("E24").Value = '  +3.61%  '

$ws.Range("D25").Value = '97.25'
$ws.Range("E25").Value = '  -2.62%  '

$ws.Range("D26").Value = '3.99'
$ws.Range("E26").Value = '  -0.32%  '

$ws.Range("E27").Value = '  +1.76%  '

$ws.Range("E28").Value = '  -2.74%  '

$ws.Range("D29").Value = '9.48'
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("D30").Value = '32.66'
$ws.Range("E30").Value = '  +4.96%  '

$ws.Range("D31").Value = '8.51'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").Value = '6.80'
$ws.Range("E32").Value = '  +4.52%  '

$ws.Range("D33").Value = '600.40'
$ws.Range("E33").Value = '  +7.26%  '

$ws.Range("D34").Value = '3.940.40'
$ws.Range("E34").Value = '  +2.83%  '

$ws.Range("D35").Value = '10.92'
$ws.Range("E35").Value = '  +0.31%  '

$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("D37").Value = '3.44'
$ws.Range("E37").Value = '  -8.52%  '

$ws.Range("E38").Value = '  -0.28%  '

$ws.Range("D39").Value = '55.49'
$ws.Range("E39").Value = '  -0.99%  '

$ws.Range("E40").Value = '  +3.36%  '

$ws.Range("D41").Value = '3.28'
$ws.Range("E41").Value = '  +3.91%  '

$ws.Range("D42").Value = '2.66'
$ws.Range("E42").Value = '  +2.29%  '

$ws.Range("D43").Value = '32.68'
$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '3.39'
$ws.Range("E44").Value = '  +0.75%  '

$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0684'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").Value = '0.334'
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").Value = '0.0415'
$ws.Range("E47").Value = '  +1.54%  '

$ws.Range("D48").Value = '0.128'
$ws.Range("E48").Value = '  +1.50%  '

$ws.Range("E49").Value = '  +0.63%  '

$ws.Range("E50").Value = '  +8.40%  '

$ws.Range("D51").Value = '2.53'
$ws.Range("E51").Value = '  +0.41%  '
